$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 405; this shifts existing rows 405-531 down to 406-532
$ws.Rows.Item(405).Insert()

# Populate the new row 405 with its values
$ws.Cells.Item(405, 1).Value = 9
$ws.Cells.Item(405, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(405, 3).Value = 'Metropolitana'
$ws.Cells.Item(405, 4).Value = 44876
$ws.Cells.Item(405, 5).Value = 13
$ws.Cells.Item(405, 6).Value = 'Fruta'
$ws.Cells.Item(405, 7).Value = 100108
$ws.Cells.Item(405, 8).Value = 'Tropicales y subtropicales'
$ws.Cells.Item(405, 9).Value = 100108002
$ws.Cells.Item(405, 10).Value = 'Mango'
$ws.Cells.Item(405, 11).Value = 'Sin especificar'
$ws.Cells.Item(405, 12).Value = 'Primera'
$ws.Cells.Item(405, 13).Value = 620
$ws.Cells.Item(405, 14).Value = 7000
$ws.Cells.Item(405, 15).Value = 8000
$ws.Cells.Item(405, 16).Value = 7452
$ws.Cells.Item(405, 17).Value = '$/bandeja 4 kilos'
$ws.Cells.Item(405, 18).Value = 'Brasil'
$ws.Cells.Item(405, 19).Value = 1863
$ws.Cells.Item(405, 20).Value = 4
